$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.012.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.666.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.51%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.49%  "

$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.703"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.42%  "

$ws.Range("E11").Value = "  -6.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000272"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.255.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.670.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.10%  "

$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "67.800.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "399.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "87.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.11%  "

$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "67.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "45.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.83%  "

$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "610.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.394"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "

$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0759"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -15.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.134"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0425"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.809.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.135"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.87%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.11%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.77%  "
